# Applies the per-row crypto-price refresh captured in the commit diff.
#
# Quirk: several updated Price-column values are numeric-looking text
# (e.g. '2.10', '59.00') that must stay plain text (matching the authored
# inlineStr cells) instead of being auto-coerced into numbers (which would
# silently drop trailing zeros, e.g. '2.10' -> 2.1). For those cells we
# temporarily force a Text number format, assign the literal string, then
# reset the cell's style back to Normal so no stray formatting lingers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.340.17'
$ws.Range('E2').Value = '  -0.33%  '

# Row 3
$ws.Range('D3').Value = '2.175.78'
$ws.Range('E3').Value = '  -1.74%  '

# Row 4
$ws.Range('E4').Value = '  -0.08%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '252.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.91%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.612'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.04%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.14'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.31%  '

# Row 8
$ws.Range('E8').Value = '  +0.00%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.580'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.21%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.80%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0911'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.65%  '

# Row 12
$ws.Range('E12').Value = '  -0.40%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.74'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.25%  '

# Row 14
$ws.Range('D14').Value = '2.505.50'
$ws.Range('E14').Value = '  -1.81%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.11'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.56%  '

# Row 16
$ws.Range('D16').Value = '2.165.92'
$ws.Range('E16').Value = '  -2.29%  '

# Row 17
$ws.Range('E17').Value = '  -4.05%  '

# Row 18
$ws.Range('D18').Value = '42.269.08'
$ws.Range('E18').Value = '  -0.27%  '

# Row 19
$ws.Range('E19').Value = '  -3.39%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.49'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.22%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.09%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.30%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.36'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.74%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.10'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.45%  '

# Row 25
$ws.Range('E25').Value = '  -0.12%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.42'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.39%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.89%  '

# Row 28
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.22'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.01%  '

# Row 29
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.17'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.06%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.69'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.02%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.51'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.05%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.96'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.38%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0807'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.02%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.07'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.59%  '

# Row 35
$ws.Range('E35').Value = '  -1.55%  '

# Row 36
$ws.Range('E36').Value = '  -0.86%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.18'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.08%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0334'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.25%  '

# Row 39
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.04'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.49%  '

# Row 40
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '11.68'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.10%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.195'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.05%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '59.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.45%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.10'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.95%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.99'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.76%  '

# Row 45
$ws.Range('E45').Value = '  +8.19%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0971'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.85%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.459'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +8.90%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.51%  '

# Row 49
$ws.Range('E49').Value = '  -2.07%  '

# Row 50
$ws.Range('E50').Value = '  -0.84%  '

# Row 51
$ws.Range('E51').Value = '  +0.11%  '
